# Applies the 2023-02-17 crypto symbol-list refresh described by the diff:
# - Updates Price (D) and Volume(1h) (E) figures for most rows
# - Swaps rows 11/12 (MandalaExchangeToken <-> BitrueCoin) back to their new order
# All target cells are stored as text in the workbook, so we force a text
# number format while writing the value (then restore the original cell style)
# to avoid Excel auto-converting these strings into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Ref, $Text) {
    $range = $Worksheet.Range($Ref)
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = $savedStyle
}

Set-TextValue $ws "D2" "308.66"
Set-TextValue $ws "E2" "-4.45%"
Set-TextValue $ws "D3" "48.68"
Set-TextValue $ws "E3" "-1.66%"
Set-TextValue $ws "D4" "5.187"
Set-TextValue $ws "E4" "-3.06%"
Set-TextValue $ws "D5" "0.07748"
Set-TextValue $ws "E5" "-4.97%"
Set-TextValue $ws "D6" "4.497"
Set-TextValue $ws "E6" "-2.38%"
Set-TextValue $ws "D7" "1.336"
Set-TextValue $ws "E7" "14.93%"
Set-TextValue $ws "D8" "1.563"
Set-TextValue $ws "E8" "-6.77%"
Set-TextValue $ws "D9" "0.1222"
Set-TextValue $ws "E9" "-9.57%"
Set-TextValue $ws "D10" "0.1937"
Set-TextValue $ws "E10" "-1.73%"
Set-TextValue $ws "B11" "BitrueCoin"
Set-TextValue $ws "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D11" "0.04649"
Set-TextValue $ws "E11" "1.73%"
Set-TextValue $ws "B12" "MandalaExchangeToken"
Set-TextValue $ws "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D12" "0.09265"
Set-TextValue $ws "E12" "-3.00%"
Set-TextValue $ws "D13" "0.1047"
Set-TextValue $ws "E13" "0.00%"
Set-TextValue $ws "D14" "0.001259"
Set-TextValue $ws "E14" "-5.86%"
Set-TextValue $ws "D15" "0.04177"
Set-TextValue $ws "E15" "-2.82%"
Set-TextValue $ws "D16" "0.005852"
Set-TextValue $ws "E16" "-0.04%"
Set-TextValue $ws "D17" "3.327"
Set-TextValue $ws "E17" "-1.83%"
Set-TextValue $ws "D18" "2.272"
Set-TextValue $ws "E18" "-6.85%"
Set-TextValue $ws "E19" "2.72%"
Set-TextValue $ws "D20" "8.356"
Set-TextValue $ws "E20" "2.95%"
Set-TextValue $ws "D21" "0.1340"
Set-TextValue $ws "E21" "-4.91%"
Set-TextValue $ws "D23" "0.001275"
Set-TextValue $ws "E23" "-2.34%"
Set-TextValue $ws "D24" "0.004170"
Set-TextValue $ws "E24" "-3.21%"
Set-TextValue $ws "D25" "0.0001351"
Set-TextValue $ws "E25" "0.06%"
Set-TextValue $ws "E26" "-4.03%"
Set-TextValue $ws "D38" "0.02556"
Set-TextValue $ws "E38" "-7.74%"
Set-TextValue $ws "D39" "0.05851"
Set-TextValue $ws "E39" "5.81%"
Set-TextValue $ws "D40" "0.01075"
Set-TextValue $ws "E40" "73.38%"
Set-TextValue $ws "D41" "0.007925"
Set-TextValue $ws "E41" "2.25%"
Set-TextValue $ws "D42" "0.1420"
Set-TextValue $ws "E42" "-1.93%"
Set-TextValue $ws "D43" "0.008357"
Set-TextValue $ws "E43" "8.77%"
Set-TextValue $ws "D44" "0.007699"
Set-TextValue $ws "E44" "-4.92%"
Set-TextValue $ws "D45" "0.3090"
Set-TextValue $ws "E45" "-11.92%"
Set-TextValue $ws "D46" "0.00006966"
Set-TextValue $ws "E46" "2.90%"
Set-TextValue $ws "D47" "0.00000000752"
Set-TextValue $ws "E47" "0.23%"
Set-TextValue $ws "D48" "0.05670"
Set-TextValue $ws "E48" "-7.54%"
Set-TextValue $ws "E49" "0.16%"
Set-TextValue $ws "D50" "0.00002104"
Set-TextValue $ws "E50" "0.23%"
Set-TextValue $ws "D51" "0.0002004"
Set-TextValue $ws "E51" "0.23%"
